$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @{ Row=2; B='Bitcoin'; C='https://coinranking.com/coin/Qwsogvtv82FCd+bitcoin-btc'; D='29.556.22'; E='  +2.29%  ' },
    @{ Row=3; B='Ethereum'; C='https://coinranking.com/coin/razxDUgYGNAdQ+ethereum-eth'; D='1.996.33'; E='  +6.18%  ' },
    @{ Row=4; B='TetherUSD'; C='https://coinranking.com/coin/HIVsRcGKkPFtW+tetherusd-usdt'; D='1.001'; E='  -0.02%  ' },
    @{ Row=5; B='BNB'; C='https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb'; D='327.35'; E='  +0.85%  ' },
    @{ Row=6; B='USDC'; C='https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc'; D='1.001'; E='  -0.05%  ' },
    @{ Row=7; B='XRP'; C='https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp'; D='0.4686'; E='  +1.37%  ' },
    @{ Row=8; B='Cardano'; C='https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'; D='0.3946'; E='  +1.64%  ' },
    @{ Row=9; B='OKB'; C='https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'; D='46.64'; E='  -0.18%  ' },
    @{ Row=10; B='Dogecoin'; C='https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'; D='0.08038'; E='  +2.32%  ' },
    @{ Row=11; B='Polygon'; C='https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'; D='1.001'; E='  +1.63%  ' },
    @{ Row=12; B='Solana'; C='https://coinranking.com/coin/zNZHO_Sjf+solana-sol'; D='22.87'; E='  +5.02%  ' },
    @{ Row=13; B='WrappedEther'; C='https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'; D='1.994.94'; E='  +7.47%  ' },
    @{ Row=14; B='Chainlink'; C='https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'; D='7.242'; E='  +3.40%  ' },
    @{ Row=15; B='Polkadot'; C='https://coinranking.com/coin/25W7FG7om+polkadot-dot'; D='5.863'; E='  +3.32%  ' },
    @{ Row=16; B='TRON'; C='https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'; D='0.07130'; E='  +2.38%  ' },
    @{ Row=17; B='Litecoin'; C='https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'; D='89.04'; E='  +0.42%  ' },
    @{ Row=18; B='BinanceUSD'; C='https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'; D='1.004'; E='  +0.03%  ' },
    @{ Row=19; B='ShibaInu'; C='https://coinranking.com/coin/xz24e0BjL+shibainu-shib'; D='0.00001004'; E='  +1.09%  ' },
    @{ Row=20; B='Avalanche'; C='https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'; D='17.43'; E='  +2.76%  ' },
    @{ Row=21; B='Dai'; C='https://coinranking.com/coin/MoTuySvg7+dai-dai'; D='1.000'; E='  -0.11%  ' },
    @{ Row=22; B='WrappedBTC'; C='https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'; D='29.578.82'; E='  +2.36%  ' },
    @{ Row=23; B='Uniswap'; C='https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'; D='5.569'; E='  +5.60%  ' },
    @{ Row=24; B='Cosmos'; C='https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'; D='11.25'; E='  +2.22%  ' },
    @{ Row=25; B='Toncoin'; C='https://coinranking.com/coin/67YlI0K1b+toncoin-ton'; D='2.102'; E='  -0.02%  ' },
    @{ Row=26; B='Monero'; C='https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'; D='158.12'; E='  +1.21%  ' },
    @{ Row=27; B='EthereumClassic'; C='https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'; D='19.72'; E='  +2.05%  ' },
    @{ Row=28; B='InternetComputer(DFINITY)'; C='https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'; D='5.969'; E='  +1.06%  ' },
    @{ Row=29; B='BitcoinCash'; C='https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'; D='120.40'; E='  +2.22%  ' },
    @{ Row=30; B='LidoDAOToken'; C='https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'; D='1.945'; E='  +2.19%  ' },
    @{ Row=31; B='Stellar'; C='https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'; D='0.09457'; E='  +0.96%  ' },
    @{ Row=32; B='ImmutableX'; C='https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'; D='0.9174'; E='  +1.84%  ' },
    @{ Row=33; B='ARBITRUM'; C='https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'; D='1.363'; E='  +3.19%  ' },
    @{ Row=34; B='Filecoin'; C='https://coinranking.com/coin/ymQub4fuB+filecoin-fil'; D='5.285'; E='  +0.27%  ' },
    @{ Row=35; B='HuobiToken'; C='https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'; D='3.219'; E='  -0.94%  ' },
    @{ Row=36; B='PEPE'; C='https://coinranking.com/coin/03WI8NQPF+pepe-pepe'; D='0.000003509'; E='  +82.32%  ' },
    @{ Row=37; B='Hedera'; C='https://coinranking.com/coin/jad286TjB+hedera-hbar'; D='0.05837'; E='  +1.46%  ' },
    @{ Row=38; B='TrustWalletToken'; C='https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'; D='1.178'; E='  +0.52%  ' },
    @{ Row=39; B='VeChain'; C='https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'; D='0.02124'; E='  +2.14%  ' },
    @{ Row=40; B='FraxShare'; C='https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'; D='7.907'; E='  +3.54%  ' },
    @{ Row=41; B='TheSandbox'; C='https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'; D='0.5774'; E='  +1.96%  ' },
    @{ Row=42; B='Algorand'; C='https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'; D='0.1825'; E='  +2.77%  ' },
    @{ Row=43; B='Aptos'; C='https://coinranking.com/coin/HGYj5JCv5+aptos-apt'; D='9.892'; E='  +2.05%  ' },
    @{ Row=44; B='MXToken'; C='https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'; D='2.805'; E='  +10.29%  ' },
    @{ Row=45; B='EnergySwap'; C='https://coinranking.com/coin/SbWqqTui-+energyswap-ens'; D='12.08'; E='  +0.80%  ' },
    @{ Row=46; B='Decentraland'; C='https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'; D='0.5405'; E='  +1.23%  ' },
    @{ Row=47; B='RenderToken'; C='https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'; D='2.217'; E='  -0.87%  ' },
    @{ Row=48; B='Cronos'; C='https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'; D='0.06979'; E='  -0.88%  ' },
    @{ Row=49; B='NEARProtocol'; C='https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'; D='1.870'; E='  +1.31%  ' },
    @{ Row=50; B='Quant'; C='https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'; D='114.11'; E='  +1.49%  ' },
    @{ Row=51; B='WOONetwork'; C='https://coinranking.com/coin/k-J3YwacF+woonetwork-woo'; D='0.3081'; E='  +8.55%  ' }
)

# Force column D to be read as text so numeric-looking strings
# like "29.556.22" or "1.001" are preserved verbatim.
$dRange = $ws.Range("D2:D51")
$dRange.NumberFormat = "@"

foreach ($r in $rows) {
    $ws.Range("B" + $r.Row).Value = $r.B
    $ws.Range("C" + $r.Row).Value = $r.C
    $ws.Range("D" + $r.Row).Value = $r.D
    $ws.Range("E" + $r.Row).Value = $r.E
}

# Restore the default (unstyled) cell style now that values are text.
$dRange.Style = "Normal"
